# Update countries & provincias Spain
# Applies updated case-count data and resulting rank-order swaps
# from the "Pais" sheet, matching the upstream data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'Datos actualizados a 5 de Agosto de 2020 a las 17:04'
$ws.Range("B4").Value = 4920917
$ws.Range("C4").Value = 2497
$ws.Range("D4").Value = 2483162
$ws.Range("E4").Value = 2277383
$ws.Range("G4").Value = 82
$ws.Range("H4").Value = 160372
$ws.Range("B6").Value = 1926642
$ws.Range("C6").Value = 20029
$ws.Range("D6").Value = 1298528
$ws.Range("E6").Value = 588033
$ws.Range("G6").Value = 261
$ws.Range("H6").Value = 40081
$ws.Range("B15").Value = 307184
$ws.Range("C15").Value = 891
$ws.Range("D21").Value = 96948
$ws.Range("E21").Value = 112578
$ws.Range("G21").Value = 30
$ws.Range("H21").Value = 4009
$ws.Range("B22").Value = 213423
$ws.Range("C22").Value = 343
$ws.Range("E22").Value = 9483
$ws.Range("B25").Value = 117878
$ws.Range("C25").Value = 86
$ws.Range("D25").Value = 102596
$ws.Range("E25").Value = 6324
$ws.Range("A37").Value = 'Republica Dominicana'
$ws.Range("B37").Value = 75660
$ws.Range("C37").Value = 1365
$ws.Range("D37").Value = 40122
$ws.Range("E37").Value = 34316
$ws.Range("G37").Value = 9
$ws.Range("H37").Value = 1222
$ws.Range("A38").Value = 'Ucrania'
$ws.Range("B38").Value = 75490
$ws.Range("C38").Value = 1271
$ws.Range("D38").Value = 41527
$ws.Range("E38").Value = 32175
$ws.Range("G38").Value = 24
$ws.Range("H38").Value = 1788
$ws.Range("D46").Value = 47768
$ws.Range("E46").Value = 6459
$ws.Range("B47").Value = 53509
$ws.Range("C47").Value = 1144
$ws.Range("D47").Value = 41199
$ws.Range("E47").Value = 10238
$ws.Range("G47").Value = 35
$ws.Range("H47").Value = 2072
$ws.Range("B65").Value = 26222
$ws.Range("C65").Value = 408
$ws.Range("E65").Value = 7232
$ws.Range("G65").Value = 13
$ws.Range("H65").Value = 823
$ws.Range("B66").Value = 23873
$ws.Range("C66").Value = 671
$ws.Range("D66").Value = 9930
$ws.Range("E66").Value = 13552
$ws.Range("G66").Value = 3
$ws.Range("H66").Value = 391
$ws.Range("B70").Value = 20336
$ws.Range("C70").Value = 459
$ws.Range("D70").Value = 8598
$ws.Range("E70").Value = 11382
$ws.Range("G70").Value = 13
$ws.Range("H70").Value = 356
$ws.Range("B91").Value = 7625
$ws.Range("C91").Value = 42
$ws.Range("D91").Value = 6399
$ws.Range("E91").Value = 1165
$ws.Range("A95").Value = 'Zambia'
$ws.Range("B95").Value = 7022
$ws.Range("C95").Value = 229
$ws.Range("D95").Value = 5667
$ws.Range("E95").Value = 1179
$ws.Range("G95").Value = 3
$ws.Range("H95").Value = 176
$ws.Range("A96").Value = 'Luxemburgo'
$ws.Range("B96").Value = 6917
$ws.Range("D96").Value = 5537
$ws.Range("E96").Value = 1262
$ws.Range("H96").Value = 118
$ws.Range("A101").Value = 'Republica de Yibuti'
$ws.Range("B101").Value = 5330
$ws.Range("C101").Value = 82
$ws.Range("D101").Value = 5057
$ws.Range("E101").Value = 214
$ws.Range("H101").Value = 59
$ws.Range("A102").Value = 'Libano'
$ws.Range("B102").Value = 5271
$ws.Range("D102").Value = 1837
$ws.Range("E102").Value = 3369
$ws.Range("H102").Value = 65
$ws.Range("B119").Value = 2838
$ws.Range("C119").Value = 4
$ws.Range("E119").Value = 290
$ws.Range("B123").Value = 2540
$ws.Range("C123").Value = 70
$ws.Range("D123").Value = 556
$ws.Range("E123").Value = 1972
$ws.Range("B143").Value = 1221
$ws.Range("C143").Value = 5
$ws.Range("D143").Value = 699
$ws.Range("E143").Value = 444
$ws.Range("B144").Value = 1213
$ws.Range("C144").Value = 10
$ws.Range("D144").Value = 1102
$ws.Range("E144").Value = 106
$ws.Range("A157").Value = 'Gambia'
$ws.Range("B157").Value = 799
$ws.Range("C157").Value = 128
$ws.Range("D157").Value = 115
$ws.Range("E157").Value = 668
$ws.Range("G157").Value = 2
$ws.Range("H157").Value = 16
$ws.Range("A158").Value = 'Lesoto'
$ws.Range("B158").Value = 726
$ws.Range("D158").Value = 174
$ws.Range("E158").Value = 531
$ws.Range("H158").Value = 21
$ws.Range("A159").Value = 'Bahamas'
$ws.Range("B159").Value = 715
$ws.Range("C159").Value = 0
$ws.Range("D159").Value = 91
$ws.Range("E159").Value = 610
$ws.Range("H159").Value = 14
$ws.Range("A160").Value = 'Vietnam'
$ws.Range("B160").Value = 713
$ws.Range("C160").Value = 41
$ws.Range("D160").Value = 381
$ws.Range("E160").Value = 324
$ws.Range("H160").Value = 8
$ws.Range("A161").Value = 'Crucero'
$ws.Range("B161").Value = 712
$ws.Range("D161").Value = 651
$ws.Range("E161").Value = 48
$ws.Range("H161").Value = 13
$ws.Range("A162").Value = 'San Marino'
$ws.Range("B162").Value = 699
$ws.Range("D162").Value = 657
$ws.Range("E162").Value = 0
$ws.Range("H162").Value = 42
$ws.Range("B169").Value = 357
$ws.Range("C169").Value = 2
$ws.Range("D169").Value = 305
$ws.Range("E169").Value = 46
$ws.Range("B179").Value = 197
$ws.Range("C179").Value = 3
$ws.Range("E179").Value = 54
